$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
}

# Paragraph 1: Title "LTI" -> "LTI - Java"
Replace-InParagraph 1 "LTI" "LTI - Java"

# Paragraph 2: Duration "2 hours" -> "Duration : 1 hours"
Replace-InParagraph 2 "2 hours" "Duration : 1 hours"

# Paragraph 4: Question 1
Replace-InParagraph 4 "1)what is 2+5?" "1)adf"
Replace-InParagraph 4 "a)3" "a)d"
Replace-InParagraph 4 "b)7" "b)s"
Replace-InParagraph 4 "c)5" "c)d"
Replace-InParagraph 4 "d)8" "d)w"

# Paragraph 5: Question 2
Replace-InParagraph 5 "2)2" "2)ffff"
Replace-InParagraph 5 "a)3" "a)1"
Replace-InParagraph 5 "c)1" "c)3"
Replace-InParagraph 5 "d)3" "d)4"

# Paragraph 6: Question 3
Replace-InParagraph 6 "3)what ?" "3)dfa"
Replace-InParagraph 6 "a)yes" "a)d"
Replace-InParagraph 6 "b)no" "b)s"
Replace-InParagraph 6 "c)o" "c)f"
Replace-InParagraph 6 "d)kk" "d)s"
